# Generate Report for Handback
# Rename the handoff/handback source file identifiers and refresh the
# handoff/handback timestamps + generated .xlf file names to reflect the
# new handback report run.

$wb = $excel.ActiveWorkbook

# ---- Overview sheet ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = "768f0ad2-e387-44ed-864f-24caa32fb264.md"
$wsOverview.Range("A3").Value = "ffff75c99d35-d159-48af-8e87-3e5499d2c770.md"

# ---- zh-cn sheet ----
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A2").Value = "768f0ad2-e387-44ed-864f-24caa32fb264.md"
$wsZh.Range("C2").Value = "768f0ad2-e387-44ed-864f-24caa32fb264.faeb1841e910f0f247a6b58eb32b1300dfadfe6b.zh-cn.xlf"
$wsZh.Range("D2").Value = "2016-03-09 03:43:48"
$wsZh.Range("E2").Value = "768f0ad2-e387-44ed-864f-24caa32fb264.md"
$wsZh.Range("F2").Value = "768f0ad2-e387-44ed-864f-24caa32fb264.faeb1841e910f0f247a6b58eb32b1300dfadfe6b.zh-cn.xlf"
$wsZh.Range("G2").Value = "2016-03-09 03:44:37"

$wsZh.Range("A3").Value = "ffff75c99d35-d159-48af-8e87-3e5499d2c770.md"
$wsZh.Range("C3").Value = "768f0ad2-e387-44ed-864f-24caa32fb264.faeb1841e910f0f247a6b58eb32b1300dfadfe6b.zh-cn.xlf"
$wsZh.Range("D3").Value = "2016-03-09 03:43:48"
$wsZh.Range("E3").Value = "768f0ad2-e387-44ed-864f-24caa32fb264.md"
$wsZh.Range("F3").Value = "768f0ad2-e387-44ed-864f-24caa32fb264.faeb1841e910f0f247a6b58eb32b1300dfadfe6b.zh-cn.xlf"
$wsZh.Range("G3").Value = "2016-03-09 03:44:37"

# ---- de-de sheet ----
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A2").Value = "768f0ad2-e387-44ed-864f-24caa32fb264.md"
$wsDe.Range("C2").Value = "768f0ad2-e387-44ed-864f-24caa32fb264.faeb1841e910f0f247a6b58eb32b1300dfadfe6b.de-de.xlf"
$wsDe.Range("D2").Value = "2016-03-09 03:43:58"
$wsDe.Range("E2").Value = "768f0ad2-e387-44ed-864f-24caa32fb264.md"
$wsDe.Range("F2").Value = "768f0ad2-e387-44ed-864f-24caa32fb264.faeb1841e910f0f247a6b58eb32b1300dfadfe6b.de-de.xlf"
$wsDe.Range("G2").Value = "2016-03-09 03:45:14"

$wsDe.Range("A3").Value = "ffff75c99d35-d159-48af-8e87-3e5499d2c770.md"
$wsDe.Range("C3").Value = "768f0ad2-e387-44ed-864f-24caa32fb264.faeb1841e910f0f247a6b58eb32b1300dfadfe6b.de-de.xlf"
$wsDe.Range("D3").Value = "2016-03-09 03:43:58"
$wsDe.Range("E3").Value = "768f0ad2-e387-44ed-864f-24caa32fb264.md"
$wsDe.Range("F3").Value = "768f0ad2-e387-44ed-864f-24caa32fb264.faeb1841e910f0f247a6b58eb32b1300dfadfe6b.de-de.xlf"
$wsDe.Range("G3").Value = "2016-03-09 03:45:14"
